$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Simple text replacements (row count unaffected) ---
$t.Rows.Item(1).Cells.Item(1).Range.Text  = "0M"
$t.Rows.Item(2).Cells.Item(1).Range.Text  = "0M"
$t.Rows.Item(3).Cells.Item(1).Range.Text  = "0M"
$t.Rows.Item(4).Cells.Item(1).Range.Text  = "269"
$t.Rows.Item(5).Cells.Item(1).Range.Text  = "0.00003"
$t.Rows.Item(6).Cells.Item(1).Range.Text  = "0.00073"

# Collapse the three multi-run summary rows (near the end of the table)
# down to a single value each.
$t.Rows.Item(44).Cells.Item(1).Range.Text = "100"
$t.Rows.Item(45).Cells.Item(1).Range.Text = "0.06"
$t.Rows.Item(46).Cells.Item(1).Range.Text = "2430"

# --- Remove the four rows that followed the old "0.00042" row ---
# (was rows 7,8,9,10 holding 0.00013 / 0.00004 / 0.00014 / 0.00017)
$t.Rows.Item(7).Delete()
$t.Rows.Item(7).Delete()
$t.Rows.Item(7).Delete()
$t.Rows.Item(7).Delete()

# After the four deletions, the old row 12 (0.01667) is now row 8.
$t.Rows.Item(8).Cells.Item(1).Range.Text = "0.00007"

# --- Insert four new rows after the row that used to hold 0.01667 ---
# Each Add(beforeRow) call inserts immediately above the same anchor row,
# so rows must be added in reverse order to end up in forward order.
$beforeRow = $t.Rows.Item(9)

$newRow4 = $t.Rows.Add($beforeRow)
$newRow4.Cells.Item(1).Range.Text = "0.06099"

$newRow3 = $t.Rows.Add($beforeRow)
$newRow3.Cells.Item(1).Range.Text = "0.00050"

$newRow2 = $t.Rows.Add($beforeRow)
$newRow2.Cells.Item(1).Range.Text = "0.00034"

$newRow1 = $t.Rows.Add($beforeRow)
$newRow1.Cells.Item(1).Range.Text = "0.00027"
